$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(16, 1).Value = "'11/28/2025"
$ws.Cells.Item(16, 1).ClearFormats()
$ws.Cells.Item(16, 2).Value = 418.3989999999994
$ws.Cells.Item(16, 3).Value = 0.05915406107567187
$ws.Cells.Item(16, 4).Value = 25
